$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'thermal men leggings'
    2 = 'compression pants with knee pads'
    3 = 'under armor compression pants men'
    4 = 'compression pants with pads'
    5 = 'basketball warm up pants men'
    6 = 'basketball compression leggings'
    7 = 'bjj compression pants'
    8 = 'padded leggings basketball'
    9 = 'football tights'
    10 = 'knee pads for volleyball youth'
    11 = 'recovery compression pants men'
    12 = 'hockey compression leggings'
    13 = 'knee compressions'
    14 = 'men sports tights'
    15 = 'wrestling knee pads men'
    16 = 'mens running leggings'
    17 = 'thick knee pad'
    18 = 'football pads pants'
    19 = 'boys athletic pants'
    20 = 'running compression tights'
    21 = 'volleyball gear men'
    22 = 'girls capri leggings'
    23 = 'black capri leggings'
    24 = 'coolomg basketball knee pads'
    25 = 'tight with knee pads'
    26 = 'tights mens nike'
    27 = 'muscle compression pants for men'
    28 = 'hayabusa compression pants men'
    29 = 'girls athletic leggings'
    30 = 'black capri leggings for women'
    31 = 'asics youth knee pads'
    32 = 'men''s tights sports'
    33 = 'tights with knee pads'
    34 = 'lavento compression pants'
    35 = 'tsla compression pants men'
    36 = 'mens black baseball pants'
    37 = 'knee pad protector'
    38 = 'compression pants boys basketball'
    39 = 'black leggings men'
    40 = 'athletic leggings men'
    41 = 'youth hex knee pads'
    42 = 'wrestling youth knee pads'
    43 = 'adult volleyball knee pads'
    44 = 'basketball shorts with pads'
    45 = 'mens wrestling tights'
    46 = 'softball sliding pants'
    47 = 'black baseball pants'
    48 = 'wrestling tights for men'
    49 = 'mens running capris'
    50 = 'mens football pants with pads'
    51 = 'soccer pants'
    52 = 'baseball sliding shorts boys'
    53 = 'compression knee guards'
    54 = 'knees protector'
    55 = 'baseball pants youth'
    56 = 'mens basketball shorts black'
    57 = 'mens compression'
    58 = 'wrestling clothes for men'
    59 = 'mens sports pants'
    60 = 'boys gym pants'
    61 = 'knee pads football adult'
    62 = 'leggings men short'
    63 = 'running pants men tall'
    64 = 'boys running tights youth'
    65 = 'softball pants for girls youth'
    66 = 'little boys athletic pants'
    67 = 'boys running pants size'
    68 = 'football youth pants'
    69 = 'comfortable knee pads'
    70 = 'knee sleeve with pad'
    71 = 'women compression tights'
    72 = 'knee pads toddler'
    73 = 'knee pads basketball kids'
    74 = 'knee protector for kids'
    75 = 'dark green knee pads for basketball'
    76 = 'skateboarding knee pads youth'
    77 = 'compression pants with padding basketball'
    78 = 'firefighter compression pants'
    79 = 'skateboard knee and elbow pads youth'
    80 = 'skateboard youth knee pads'
    81 = 'warm up pants men'
    82 = 'mc david knee pad'
    83 = 'compression knee pads men'
    84 = 'men basketball pants'
    85 = 'premium knee pad'
    86 = 'kids compression pants for basketball'
    87 = 'toddler knee pad'
    88 = 'padded knee sleeves men'
    89 = 'mtb knee pads for men'
    90 = 'compression knee sleeves pads'
    91 = 'basketball knee pads leggings'
    92 = 'soccer compression pants men'
    93 = 'mens small leggings'
    94 = 'compressions pants mens'
    95 = 'capri basketball leggings for boys'
    96 = 'youth small black baseball pants'
    97 = 'mens black leggings'
    98 = 'basket ball knee pads youth'
    99 = 'bjj leggings men'
    100 = 'mens gym tights'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
